# "updated activity till excel form"
#
# The sheet holds David Warner's per-innings batting activity
# (runs/balls/fours/sixes in columns C:F, one innings per row).
# This edit refreshes those figures by re-pointing each row at a
# different innings' numbers (the values were reshuffled across rows
# 2-7 and 10-17; rows 8-9 and a handful of unchanged cells are left
# untouched).
#
# The source cells are numbers-stored-as-text (OOXML t="str"), so we
# mark each destination cell as Text (NumberFormat "@") before typing
# the new value - otherwise Excel's automatic type detection would
# silently convert these numeric-looking strings back into real
# numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark the touched cells as Text so the numeric-looking values stick
# as strings instead of being auto-converted to numbers.
$ws.Range("C2:F2").NumberFormat = "@"
$ws.Range("C3:D3").NumberFormat = "@"
$ws.Range("C4:E4").NumberFormat = "@"
$ws.Range("C5:F5").NumberFormat = "@"
$ws.Range("C6:F6").NumberFormat = "@"
$ws.Range("C7:E7").NumberFormat = "@"
$ws.Range("C10:F10").NumberFormat = "@"
$ws.Range("C11:F11").NumberFormat = "@"
$ws.Range("C12:D12").NumberFormat = "@"
$ws.Range("F12").NumberFormat = "@"
$ws.Range("C13:F13").NumberFormat = "@"
$ws.Range("C14:F14").NumberFormat = "@"
$ws.Range("C15:E15").NumberFormat = "@"
$ws.Range("C16:F16").NumberFormat = "@"
$ws.Range("C17:F17").NumberFormat = "@"

# Row 2
$ws.Range("C2").Value = "17"
$ws.Range("D2").Value = "17"
$ws.Range("E2").Value = "3"
$ws.Range("F2").Value = "0"

# Row 3
$ws.Range("C3").Value = "35"
$ws.Range("D3").Value = "20"

# Row 4
$ws.Range("C4").Value = "8"
$ws.Range("D4").Value = "5"
$ws.Range("E4").Value = "0"

# Row 5
$ws.Range("C5").Value = "66"
$ws.Range("D5").Value = "34"
$ws.Range("E5").Value = "8"
$ws.Range("F5").Value = "2"

# Row 6
$ws.Range("C6").Value = "4"
$ws.Range("D6").Value = "4"
$ws.Range("E6").Value = "1"
$ws.Range("F6").Value = "0"

# Row 7
$ws.Range("C7").Value = "85"
$ws.Range("D7").Value = "58"
$ws.Range("E7").Value = "10"

# Row 10
$ws.Range("C10").Value = "47"
$ws.Range("D10").Value = "33"
$ws.Range("E10").Value = "5"
$ws.Range("F10").Value = "0"

# Row 11
$ws.Range("C11").Value = "48"
$ws.Range("D11").Value = "38"
$ws.Range("E11").Value = "3"
$ws.Range("F11").Value = "2"

# Row 12
$ws.Range("C12").Value = "28"
$ws.Range("D12").Value = "29"
$ws.Range("F12").Value = "0"

# Row 13
$ws.Range("C13").Value = "36"
$ws.Range("D13").Value = "30"
$ws.Range("E13").Value = "2"
$ws.Range("F13").Value = "1"

# Row 14
$ws.Range("C14").Value = "60"
$ws.Range("D14").Value = "44"
$ws.Range("E14").Value = "5"
$ws.Range("F14").Value = "2"

# Row 15
$ws.Range("C15").Value = "45"
$ws.Range("D15").Value = "33"
$ws.Range("E15").Value = "3"

# Row 16
$ws.Range("C16").Value = "52"
$ws.Range("D16").Value = "40"
$ws.Range("E16").Value = "5"
$ws.Range("F16").Value = "1"

# Row 17
$ws.Range("C17").Value = "6"
$ws.Range("D17").Value = "6"
$ws.Range("E17").Value = "1"
$ws.Range("F17").Value = "0"
